# Brewing Potions.xlsx -- apply commit changes
#
# Summary of changes applied:
#  1. Potions!D41: 3 -> 2
#  2. Potions!L25:L32: new "<N> AC" labels next to the Barkskin rows
#  3. Potions: insert a "Cure Moderate Wounds" CL10 row (after the existing
#     CL9 row) and six more "Cure Serious Wound" rows (CL4-CL10), pushing
#     the existing CL3 / Fly / Heroism rows down.
#  4. Active sheet changes from "Brewing" to "Potions" (with a specific
#     selected cell on each).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Potions")

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

function Copy-RowFormat($srcRow, $dstRow) {
    $ws.Range("A$srcRow`:K$srcRow").Copy() | Out-Null
    $ws.Range("A$dstRow`:K$dstRow").PasteSpecial($xlPasteFormats) | Out-Null
}

# --- 1. Simple value edit -------------------------------------------------
$ws.Range("D41").Value2 = 2

# --- 2. New "AC" labels for the Barkskin block (rows 25-32) --------------
$ws.Range("L25").Value2 = "2 AC"
$ws.Range("L26").Value2 = "2 AC"
$ws.Range("L27").Value2 = "2 AC"
$ws.Range("L28").Value2 = "2 AC"
$ws.Range("L29").Value2 = "3 AC"
$ws.Range("L30").Value2 = "3 AC"
$ws.Range("L31").Value2 = "3 AC"
$ws.Range("L32").Value2 = "4 AC"

# --- 3. Relocate existing rows 49-52 to make room for new data -----------
# Old layout (rows 49-52):
#   49 blank separator row
#   50 Cure Serious Wound (CL3)
#   51 Fly
#   52 Heroism
# New layout (rows 49-60):
#   49 Cure Moderate Wounds (CL10)   <- NEW
#   50 blank separator row           <- moved from 49
#   51 Cure Serious Wound (CL3)      <- moved from 50 (D changes 5 -> 3)
#   52 Cure Serious Wound (CL4)      <- NEW
#   53 Cure Serious Wound (CL5)      <- NEW
#   54 Cure Serious Wound (CL6)      <- NEW
#   55 Cure Serious Wound (CL7)      <- NEW
#   56 Cure Serious Wound (CL8)      <- NEW
#   57 Cure Serious Wound (CL9)      <- NEW
#   58 Cure Serious Wound (CL10)     <- NEW
#   59 Fly                           <- moved from 51
#   60 Heroism                       <- moved from 52

# Move from the bottom up so we never overwrite data we still need.
Copy-RowFormat 52 60
$ws.Range("A60").Value2 = "Heroism"
$ws.Range("B60").Value2 = 3
$ws.Range("C60").Value2 = 3
$ws.Range("D60").Formula = "=10*C60"
$ws.Range("E60").Value2 = "Min"
$ws.Range("F60").Formula = "=(B60*C60*50)"
$ws.Range("G60").Formula = "=F60-(F60*0.05)"
$ws.Range("H60").Formula = "=G60/2"
$ws.Range("I60").Value2 = 8
$ws.Range("J60").Formula = "=I60/2"
$ws.Range("K60").Formula = "=5+C60"

Copy-RowFormat 51 59
$ws.Range("A59").Value2 = "Fly"
$ws.Range("B59").Value2 = 3
$ws.Range("C59").Value2 = 3
$ws.Range("D59").Formula = "=1*C59"
$ws.Range("E59").Value2 = "Min"
$ws.Range("F59").Formula = "=(B59*C59*50)"
$ws.Range("G59").Formula = "=F59-(F59*0.05)"
$ws.Range("H59").Formula = "=G59/2"
$ws.Range("I59").Value2 = 8
$ws.Range("J59").Formula = "=I59/2"
$ws.Range("K59").Formula = "=5+C59"

Copy-RowFormat 50 51
$ws.Range("A51").Value2 = "Cure Serious Wound"
$ws.Range("B51").Value2 = 3
$ws.Range("C51").Value2 = 3
$ws.Range("D51").Value2 = 3
$ws.Range("E51").Value2 = "HP"
$ws.Range("F51").Formula = "=(B51*C51*50)"
$ws.Range("G51").Formula = "=F51-(F51*0.05)"
$ws.Range("H51").Formula = "=G51/2"
$ws.Range("I51").Value2 = 8
$ws.Range("J51").Formula = "=I51/2"
$ws.Range("K51").Formula = "=5+C51"

Copy-RowFormat 49 50
$ws.Range("A50").Value2 = $null
$ws.Range("B50").Value2 = $null
$ws.Range("C50").Value2 = $null
$ws.Range("D50").Value2 = $null
$ws.Range("E50").Value2 = $null
$ws.Range("F50").Value2 = $null
$ws.Range("G50").Value2 = $null
$ws.Range("H50").Value2 = $null
$ws.Range("I50").Value2 = $null
$ws.Range("J50").Value2 = $null
$ws.Range("K50").Value2 = $null

# New row 49: Cure Moderate Wounds, CL10
Copy-RowFormat 48 49
$ws.Range("A49").Value2 = "Cure Moderate Wounds"
$ws.Range("B49").Value2 = 2
$ws.Range("C49").Value2 = 10
$ws.Range("D49").Value2 = 10
$ws.Range("E49").Value2 = "HP"
$ws.Range("F49").Formula = "=(B49*C49*50)"
$ws.Range("G49").Formula = "=F49-(F49*0.05)"
$ws.Range("H49").Formula = "=G49/2"
$ws.Range("I49").Value2 = 8
$ws.Range("J49").Formula = "=I49/2"
$ws.Range("K49").Formula = "=5+C49"

# New rows 52-58: Cure Serious Wound, CL4-CL10
$newRows = @(
    @{ Row = 52; CL = 4;  Hours = 8 },
    @{ Row = 53; CL = 5;  Hours = 8 },
    @{ Row = 54; CL = 6;  Hours = 8 },
    @{ Row = 55; CL = 7;  Hours = 16 },
    @{ Row = 56; CL = 8;  Hours = 16 },
    @{ Row = 57; CL = 9;  Hours = 16 },
    @{ Row = 58; CL = 10; Hours = 16 }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    Copy-RowFormat 51 $r
    $ws.Range("A$r").Value2 = $null
    $ws.Range("B$r").Value2 = 3
    $ws.Range("C$r").Value2 = $entry.CL
    $ws.Range("D$r").Value2 = $entry.CL
    $ws.Range("E$r").Value2 = "HP"
    $ws.Range("F$r").Formula = "=(B$r*C$r*50)"
    $ws.Range("G$r").Formula = "=F$r-(F$r*0.05)"
    $ws.Range("H$r").Formula = "=G$r/2"
    $ws.Range("I$r").Value2 = $entry.Hours
    $ws.Range("J$r").Formula = "=I$r/2"
    $ws.Range("K$r").Formula = "=5+C$r"
}

# --- 4. Switch the active sheet/selection ---------------------------------
$ws.Select()
$ws.Range("T12").Select()
